$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("template")

# Rename header E1 from "Precio_descuento" to "Descuento"
$ws.Range("E1").Value = "Descuento"

# Update discount values in column E (rows 2-260)
$discounts = @{
    2 = 10; 3 = 15; 4 = 20; 6 = 30; 7 = 35; 8 = 40; 9 = 45; 10 = 50;
    11 = 55; 12 = 60; 13 = 65; 14 = 70; 15 = 75; 16 = 80; 17 = 85; 18 = 90;
    19 = 95; 20 = 97; 21 = 99; 22 = 100
}

foreach ($row in $discounts.Keys) {
    $ws.Cells.Item($row, 5).Value = $discounts[$row]
}

for ($row = 23; $row -le 260; $row++) {
    $ws.Cells.Item($row, 5).Value = 20
}

# Update the active selection to E23
$ws.Range("E23").Select()
